$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Options")

$ws.Range("A4").Value = "Privacy Budget"
$ws.Range("B4").Value = 0.8

$ws.Range("B5").Select()
